# adding NA to lookup table and adjusting enums accordingly
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 ("0"/"Not Available"/"Not available"),
# shifting the existing rows (and the styled I6 cell) down by one.
$ws.Rows(2).Insert()

# Populate the new row 2 with the "Not applicable" lookup entry (id = -1).
# Note: there is no name_sanitized value for this row (column C stays empty).
$ws.Range("A2").Value = -1
$ws.Range("B2").Value = "Not applicable"

# Match the author's resulting selection.
$ws.Range("A3").Select()
